$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell/range in the sheet view (cosmetic change from A25 to J30)
$ws.Range("J30").Select()

# Fix wrong correlation: flip sign of column J values for rows 14 through 25 (1 -> -1)
for ($row = 14; $row -le 25; $row++) {
    $ws.Cells.Item($row, 10).Value = -1
}
